# Update faturamento (revenue) data for Bibi stores - "atualizacao dos dados da bibi"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 - Bibi Cell Mundi
$ws.Range("N2").Value = 13707.95
$ws.Range("AG2").Value = 139414.29

# Row 3 - Bibi Cell Vieiralves
$ws.Range("M3").Value = 2555
$ws.Range("N3").Value = 6272
$ws.Range("AG3").Value = 82146.8

# Row 4 - Bibi Cell Manauara
$ws.Range("M4").Value = 2668
$ws.Range("N4").Value = 3201.9
$ws.Range("AG4").Value = 43206.79

# Row 5 - Bibi Cell Ponta Negra
$ws.Range("N5").Value = 3022.01
$ws.Range("AG5").Value = 39233.13

# Row 6 - total
$ws.Range("M6").Value = 31420.44
$ws.Range("N6").Value = 26203.86
$ws.Range("AG6").Value = 304001.01
